$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "68.660.00"
Set-TextValue $ws.Range("E2") "  +0.37%  "
Set-TextValue $ws.Range("D3") "2.729.55"
Set-TextValue $ws.Range("E3") "  +3.10%  "
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "608.57"
Set-TextValue $ws.Range("E5") "  +1.79%  "
Set-TextValue $ws.Range("D6") "169.55"
Set-TextValue $ws.Range("E6") "  +6.56%  "
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("E8") "  +1.27%  "
Set-TextValue $ws.Range("D9") "2.728.69"
Set-TextValue $ws.Range("E9") "  +3.08%  "
Set-TextValue $ws.Range("D10") "0.147"
Set-TextValue $ws.Range("E10") "  +3.13%  "
Set-TextValue $ws.Range("D11") "0.367"
Set-TextValue $ws.Range("E11") "  +4.77%  "
Set-TextValue $ws.Range("E12") "  +0.98%  "
Set-TextValue $ws.Range("E13") "  -0.09%  "
Set-TextValue $ws.Range("D14") "28.76"
Set-TextValue $ws.Range("E14") "  +2.76%  "
Set-TextValue $ws.Range("D15") "3.227.56"
Set-TextValue $ws.Range("E15") "  +3.05%  "
Set-TextValue $ws.Range("E16") "  +1.93%  "
Set-TextValue $ws.Range("D17") "68.626.56"
Set-TextValue $ws.Range("E17") "  +0.44%  "
Set-TextValue $ws.Range("D18") "2.734.65"
Set-TextValue $ws.Range("E18") "  +4.17%  "
Set-TextValue $ws.Range("D19") "11.92"
Set-TextValue $ws.Range("E19") "  +4.73%  "
Set-TextValue $ws.Range("D20") "375.74"
Set-TextValue $ws.Range("E20") "  +4.33%  "
Set-TextValue $ws.Range("D21") "7.66"
Set-TextValue $ws.Range("E21") "  +3.65%  "
Set-TextValue $ws.Range("D22") "4.52"
Set-TextValue $ws.Range("E22") "  +2.16%  "
Set-TextValue $ws.Range("D23") "5.00"
Set-TextValue $ws.Range("E23") "  +4.92%  "
Set-TextValue $ws.Range("E24") "  +1.77%  "
Set-TextValue $ws.Range("D25") "73.71"
Set-TextValue $ws.Range("E25") "  -1.24%  "
Set-TextValue $ws.Range("E26") "  +0.02%  "
Set-TextValue $ws.Range("D27") "10.12"
Set-TextValue $ws.Range("E27") "  +3.64%  "
Set-TextValue $ws.Range("D28") "2.868.67"
Set-TextValue $ws.Range("E28") "  +2.98%  "
Set-TextValue $ws.Range("E29") "  +2.20%  "
Set-TextValue $ws.Range("D30") "590.66"
Set-TextValue $ws.Range("E30") "  +5.04%  "
Set-TextValue $ws.Range("E31") "  -0.14%  "
Set-TextValue $ws.Range("D32") "8.33"
Set-TextValue $ws.Range("E32") "  +3.58%  "
Set-TextValue $ws.Range("E33") "  +3.85%  "
Set-TextValue $ws.Range("D34") "1.99"
Set-TextValue $ws.Range("E34") "  +6.01%  "
Set-TextValue $ws.Range("E35") "  +3.47%  "
Set-TextValue $ws.Range("E36") "  -1.59%  "
Set-TextValue $ws.Range("E37") "  -0.03%  "
Set-TextValue $ws.Range("D38") "163.36"
Set-TextValue $ws.Range("E38") "  +3.05%  "
Set-TextValue $ws.Range("D39") "20.00"
Set-TextValue $ws.Range("E39") "  +1.61%  "
Set-TextValue $ws.Range("E40") "  +2.86%  "
Set-TextValue $ws.Range("B41") "RenderToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Range("D41") "5.49"
Set-TextValue $ws.Range("E41") "  +2.77%  "
Set-TextValue $ws.Range("B42") "Stacks"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "1.90"
Set-TextValue $ws.Range("E42") "  +1.95%  "
Set-TextValue $ws.Range("B43") "WhiteBITCoin"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D43") "17.97"
Set-TextValue $ws.Range("E43") "  +0.97%  "
Set-TextValue $ws.Range("B44") "dogwifhat"
Set-TextValue $ws.Range("C44") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D44") "2.67"
Set-TextValue $ws.Range("E44") "  +1.61%  "
Set-TextValue $ws.Range("E45") "  -0.04%  "
Set-TextValue $ws.Range("E46") "  -2.47%  "
Set-TextValue $ws.Range("D47") "40.98"
Set-TextValue $ws.Range("E47") "  +1.20%  "
Set-TextValue $ws.Range("D48") "0.604"
Set-TextValue $ws.Range("E48") "  +5.10%  "
Set-TextValue $ws.Range("D49") "155.64"
Set-TextValue $ws.Range("E49") "  -1.19%  "
Set-TextValue $ws.Range("E50") "  +3.38%  "
Set-TextValue $ws.Range("D51") "1.79"
Set-TextValue $ws.Range("E51") "  +5.84%  "
